# ---------------------------------------------------------------------------
# Adds three new worksheets (SearchProduct, ProductQuantity, AccountCreationData)
# with data/styling, and makes AccountCreationData the active sheet - matching
# the target commit "updated dataproviders and etc..".
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$missing = [System.Reflection.Missing]::Value

# Existing sheets, used as "style donors" so that new cells reuse identical
# cellXfs entries instead of creating near-duplicate styles.
$sheet1 = $wb.Worksheets.Item("Sheet1")
$billing = $wb.Worksheets.Item("Billing_Address")

# ---------------------------------------------------------------------------
# 1) Create the three new worksheets, in order, at the end of the workbook.
# ---------------------------------------------------------------------------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSearch = $wb.Worksheets.Add($missing, $lastSheet)
$wsSearch.Name = "SearchProduct"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsQty = $wb.Worksheets.Add($missing, $lastSheet)
$wsQty.Name = "ProductQuantity"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsAcct = $wb.Worksheets.Add($missing, $lastSheet)
$wsAcct.Name = "AccountCreationData"

# ---------------------------------------------------------------------------
# 2) SearchProduct sheet
# ---------------------------------------------------------------------------

$wsSearch.Range("A1").Value = "Search Product"
$wsSearch.Range("A2").Value = "Lenovo Thinkpad X1 Carbon Laptop"
$wsSearch.Range("A3").Value = "Apple MacBook Pro 13-inch"
$wsSearch.Range("A4").Value = "HP Spectre XT Pro UltraBook"
$wsSearch.Range("A5").Value = ""

# Header: yellow fill, centered, no border (same recipe as existing yellow
# header style but without the hair border around it).
$sheet1.Range("A1").Copy($wsSearch.Range("A1"))
$wsSearch.Range("A1").Value = "Search Product"
$wsSearch.Range("A1").Borders.LineStyle = -4142

# Body rows: centered, hair border (same as Billing_Address data rows).
$billing.Range("A2").Copy($wsSearch.Range("A2:A4"))
$wsSearch.Range("A2").Value = "Lenovo Thinkpad X1 Carbon Laptop"
$wsSearch.Range("A3").Value = "Apple MacBook Pro 13-inch"
$wsSearch.Range("A4").Value = "HP Spectre XT Pro UltraBook"

# Trailing blank row: blue font, wrap text, no border.
$sheet1.Range("A3").Copy($wsSearch.Range("A5"))
$wsSearch.Range("A5").Value = ""
$wsSearch.Range("A5").Borders.LineStyle = -4142
$wsSearch.Range("A5").HorizontalAlignment = 1
$wsSearch.Range("A5").WrapText = $true

$wsSearch.Columns.Item(1).ColumnWidth = 33.5

# ---------------------------------------------------------------------------
# 3) ProductQuantity sheet
# ---------------------------------------------------------------------------

$wsQty.Range("A1").Value = "Quantity"
$sheet1.Range("A1").Copy($wsQty.Range("A1"))
$wsQty.Range("A1").Value = "Quantity"

$wsQty.Range("A2").Formula = "=SUM(1)"
$wsQty.Range("A3").Formula = "=SUM(3)"
$wsQty.Range("A4").Formula = "=SUM(2)"

$billing.Range("A2").Copy($wsQty.Range("A2:A4"))
$wsQty.Range("A2").Formula = "=SUM(1)"
$wsQty.Range("A3").Formula = "=SUM(3)"
$wsQty.Range("A4").Formula = "=SUM(2)"
$wsQty.Range("A2:A4").HorizontalAlignment = -4131

# ---------------------------------------------------------------------------
# 4) AccountCreationData sheet
# ---------------------------------------------------------------------------

$headers = @("Gender", "FirstName", "LastName", "Day", "Month", "Year", "Email", "Company", "Password", "ConfirmPassword")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsAcct.Cells.Item(1, $i + 1).Value = $headers[$i]
}
$sheet1.Range("A1").Copy($wsAcct.Range("A1:J1"))
for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsAcct.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$row2 = @("male", "dinesh", "kila", 10, "June", 1998, "dineshkila99@gmail.com", "Capgemini", "dinesh123", "dinesh123")
$row3 = @("female", "nithya", "sri", 15, "August", 1995, "nithyasri@gmail.com", "HCL", "nithya123", "nithya123")
$row4 = @("male", "raju", "rajesh", 19, "April", 2001, "rajurajesh@gmail.com", "TCS", "raju1234", "raju1234")

$billing.Range("A2").Copy($wsAcct.Range("A2:J4"))

for ($i = 0; $i -lt $row2.Length; $i++) {
    $wsAcct.Cells.Item(2, $i + 1).Value = $row2[$i]
    $wsAcct.Cells.Item(3, $i + 1).Value = $row3[$i]
    $wsAcct.Cells.Item(4, $i + 1).Value = $row4[$i]
}

# Email column (G) uses the hyperlink font/style, same as Sheet1's hyperlinked
# e-mail cells.
$sheet1.Range("A3").Copy($wsAcct.Range("G2:G4"))
$wsAcct.Cells.Item(2, 7).Value = "dineshkila99@gmail.com"
$wsAcct.Cells.Item(3, 7).Value = "nithyasri@gmail.com"
$wsAcct.Cells.Item(4, 7).Value = "rajurajesh@gmail.com"

$wsAcct.Hyperlinks.Add($wsAcct.Cells.Item(2, 7), "mailto:dineshkila99@gmail.com", $missing, $missing, "dineshkila99@gmail.com")
$wsAcct.Hyperlinks.Add($wsAcct.Cells.Item(3, 7), "mailto:nithyasri@gmail.com", $missing, $missing, "nithyasri@gmail.com")
$wsAcct.Hyperlinks.Add($wsAcct.Cells.Item(4, 7), "mailto:rajurajesh@gmail.com", $missing, $missing, "rajurajesh@gmail.com")

$wsAcct.Columns.Item(7).ColumnWidth = 19.7
$wsAcct.Columns.Item(8).ColumnWidth = 12.53
$wsAcct.Columns.Item(10).ColumnWidth = 14.89

# ---------------------------------------------------------------------------
# 5) Selection / activation: AccountCreationData becomes the active sheet,
#    Billing_Address is no longer the selected tab.
# ---------------------------------------------------------------------------

$wsAcct.Range("D4").Select()
$wsAcct.Activate()

Write-Host "Workbook updated: added SearchProduct, ProductQuantity, AccountCreationData"
